$d = $word.ActiveDocument

$old = " Iquique, consta en el "
$new = " Iquique, consta en el " + "`${art8}" + "Decreto N°140/04, del Ministerio de Salud que aprobó el Reglamento Orgánico de los Servicios de Salud, "

$r = $d.Content
$found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find target text to update the Decreto del director wording."
}

Write-Host "Found:" $found
